$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.782.94'
$ws.Range('E2').Value = '  +1.06%  '

$ws.Range('D3').Value = '2.103.52'
$ws.Range('E3').Value = '  +1.99%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.96'
$ws.Range('E5').Value = '  +0.69%  '

$ws.Range('E6').Value = '  +0.69%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.07'
$ws.Range('E7').Value = '  +0.47%  '

$ws.Range('E8').Value = '  +0.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.391'
$ws.Range('E9').Value = '  +1.63%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0778'
$ws.Range('E10').Value = '  +2.49%  '

$ws.Range('E11').Value = '  +1.29%  '

$ws.Range('D12').Value = '2.416.66'
$ws.Range('E12').Value = '  +2.14%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.50'
$ws.Range('E13').Value = '  +0.77%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.29'
$ws.Range('E14').Value = '  +2.53%  '

$ws.Range('E15').Value = '  +0.93%  '

$ws.Range('E16').Value = '  +1.30%  '

$ws.Range('D17').Value = '2.117.65'
$ws.Range('E17').Value = '  +2.50%  '

$ws.Range('D18').Value = '37.775.09'
$ws.Range('E18').Value = '  +1.27%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.24'
$ws.Range('E19').Value = '  -1.27%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.25'
$ws.Range('E20').Value = '  +1.44%  '

$ws.Range('D21').Value = '0.0₃0822'
$ws.Range('E21').Value = '  +1.11%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '227.23'
$ws.Range('E22').Value = '  +0.82%  '

$ws.Range('E23').Value = '  +0.04%  '

$ws.Range('E24').Value = '  -1.54%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.41'
$ws.Range('E25').Value = '  +0.21%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.94'
$ws.Range('E26').Value = '  +0.96%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.95'
$ws.Range('E27').Value = '  +1.17%  '

$ws.Range('E28').Value = '  +2.24%  '

$ws.Range('E29').Value = '  -3.30%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.49'
$ws.Range('E30').Value = '  +1.97%  '

$ws.Range('E31').Value = '  +0.75%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.63'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.59'
$ws.Range('E33').Value = '  +0.76%  '

$ws.Range('E34').Value = '  +0.78%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.58'
$ws.Range('E35').Value = '  +0.63%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.47'
$ws.Range('E36').Value = '  +6.01%  '

$ws.Range('E37').Value = '  +1.82%  '

$ws.Range('E38').Value = '  -0.02%  '

$ws.Range('E39').Value = '  -7.68%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0999'
$ws.Range('E40').Value = '  +7.44%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.94'
$ws.Range('E41').Value = '  -0.29%  '

$ws.Range('D42').Value = '1.476.11'
$ws.Range('E42').Value = '  +0.13%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '96.68'
$ws.Range('E43').Value = '  -0.14%  '

$ws.Range('E44').Value = '  +1.26%  '

$ws.Range('E45').Value = '  +0.13%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.11'
$ws.Range('E46').Value = '  -12.22%  '

$ws.Range('E47').Value = '  +2.69%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.39'
$ws.Range('E48').Value = '  -1.12%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.32'
$ws.Range('E49').Value = '  +2.18%  '

$ws.Range('E50').Value = '  +2.97%  '

$ws.Range('D51').Value = '2.302.05'
$ws.Range('E51').Value = '  +2.19%  '
